$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old last data row (row 6 - the pre-paid cards requirement) and the
# old first data row (row 2 - the nursing/cohort requirement). Deleting row 6
# first keeps row indices for the row-2 deletion stable. This shifts the
# remaining rows (old 3,4,5) up to become rows 2,3,4, matching the target
# sheet dimension of A1:C5.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(2).Delete()

# Row 2 (index 0): Disputes System - view access
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "The Disputes System shall provide view access capability for authorized users of the application."
$ws.Cells.Item(2,3).Value = "The Disputes System shall provide view access capability for authorized users of the application, specifying the types of data or sections of the application that can be viewed and defining ""authorized users"" for clarity."

# Row 3 (index 1): Disputes System - update access
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "The Disputes System shall provide update access capability for authorized users of the application."
$ws.Cells.Item(3,3).Value = "The Disputes System shall provide update access capability for authorized users of the application, clearly defining the types of updates permitted, such as updates to user profiles and dispute details, and detailing the authentication criteria for ""authorized users."""

# Row 4 (index 2): Disputes System - select disputable transactions
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "The Disputes System must allow the users to select disputable transactions (based on the age of the transaction) from a user interface and initiate a dispute (ticket retrieval request or chargeback notification) on the selected transaction."
$ws.Cells.Item(4,3).Value = "The Disputes System must allow users to select disputable transactions based on defined criteria, including the age of the transaction, and must detail the steps for initiating a dispute, which may involve a ticket retrieval request or a chargeback notification, ensuring these actions are distinct and clearly defined."

# Row 5 (index 3, new row): Disputes System - single dispute case on multiple transactions
# Copy the formatted A4 cell into A5 so the bold/bordered/centered style carries over,
# then overwrite its value.
$ws.Cells.Item(4,1).Copy($ws.Cells.Item(5,1))
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "The Disputes System must provide the user the ability to initiate a single dispute case on multiple transactions that belong to a single merchant."
$ws.Cells.Item(5,3).Value = "The Disputes System must provide the ability for the user to initiate a dispute case that encompasses multiple transactions associated with a specific merchant, clearly defining what constitutes a ""single merchant"" and what is meant by a ""dispute case,"" including any limits on the number of transactions."
